# Edit the "Group 150" financial summary table on slide 3:
#   - resize two of the table's grid columns
#   - split the "Financials (USD MM)" header cell into two lines
#   - split each "Revenue: EBITDA:  Market Cap: Total Debt: [FTE:]" cell
#     into separate one-item-per-line paragraphs (adding an "FTE:" line
#     where it is still missing)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

# --- resize grid columns 3 and 4 (1-based) -------------------------------
$tbl.Columns.Item(3).Width = 81.99181102362205   # 1041296 EMU (was 829079)
$tbl.Columns.Item(4).Width = 354.54              # 4502658 EMU (was 4714875)

# --- header cell: "Financials (USD MM)" -> "Financials " / "(USD MM)" ---
$headerCell = $tbl.Cell(1, 3)
$headerCell.Shape.TextFrame.TextRange.Text = "Financials " + [char]13 + "(USD MM)"

# --- data rows 2..7: split the financial-summary labels into lines ------
$labelText = "Revenue:" + [char]13 + "EBITDA:" + [char]13 + "Market Cap:" + [char]13 + "Total Debt:" + [char]13 + "FTE:"
for ($r = 2; $r -le 7; $r++) {
    $cell = $tbl.Cell($r, 3)
    $cell.Shape.TextFrame.TextRange.Text = $labelText
}
